{"js": "// Neutralize racial-misclassification language to \"50M voters\" wording.\n// Three edits:\n//  1. Professional summary: \"affecting all Black and Asian-American voters,\"\n//     -> \"affecting 50M voters,\"  (plain text substitution, no run split)\n//  2. Bullet achievement: \"affecting all Black and Asian-American voters,\"\n//     -> \"affecting \" + bold/colored \"50M\" + \" voters,\"  (splits one run\n//     into three so \"50M\" can carry its own bold + color formatting)\n//  3. Impact statement: \"affecting all Black and Asian-American voters,\"\n//     -> \"affecting 50M voters nationwide,\"  (plain text substitution)\n\n// --- Edit 1: Professional summary sentence ---------------------------------\nconst summaryHits = context.document.body.search(\n  \"affecting all Black and Asian-American voters, developed geospatial ML algorithms\",\n  { matchCase: true }\n);\nsummaryHits.load(\"items\");\nawait context.sync();\n\nif (summaryHits.items.length > 0) {\n  summaryHits.items[0].insertText(\n    \"affecting 50M voters, developed geospatial ML algorithms\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n\n// --- Edit 2: Bulleted achievement (needs a bold/colored \"50M\" run) ---------\nconst bulletHits = context.document.body.search(\n  \"affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving demographic classification accuracy from\",\n  { matchCase: true }\n);\nbulletHits.load(\"items\");\nawait context.sync();\n\nif (bulletHits.items.length > 0) {\n  // Replace the whole span first with plain (unformatted) text containing\n  // the neutral \"50M\" wording, then re-find just the \"50M\" token so it can\n  // be given its own run with bold + the existing accent color.\n  bulletHits.items[0].insertText(\n    \"affecting 50M voters, developed geospatial machine learning algorithms improving demographic classification accuracy from\",\n    \"Replace\"\n  );\n  await context.sync();\n\n  const fiftyMHits = context.document.body.search(\"50M voters, developed geospatial machine learning\", {\n    matchCase: true,\n  });\n  fiftyMHits.load(\"items\");\n  await context.sync();\n\n  if (fiftyMHits.items.length > 0) {\n    // Narrow the match down to exactly \"50M\" within the hit range by\n    // re-searching scoped to that range.\n    const scoped = fiftyMHits.items[0].search(\"50M\", { matchCase: true });\n    scoped.load(\"items\");\n    await context.sync();\n\n    if (scoped.items.length > 0) {\n      const fiftyM = scoped.items[0];\n      fiftyM.font.bold = true;\n      fiftyM.font.color = \"#2C3E50\";\n      await context.sync();\n    }\n  }\n}\n\n// --- Edit 3: Impact statement sentence --------------------------------------\nconst impactHits = context.document.body.search(\n  \"affecting all Black and Asian-American voters, improved electoral prediction accuracy by 22%\",\n  { matchCase: true }\n);\nimpactHits.load(\"items\");\nawait context.sync();\n\nif (impactHits.items.length > 0) {\n  impactHits.items[0].insertText(\n    \"affecting 50M voters nationwide, improved electoral prediction accuracy by 22%\",\n    \"Replace\"\n  );\n  await context.sync();\n}\n", "ps1": "# Neutralize racial-misclassification language to \"50M voters\" wording.\n# Three edits, each located by scanning paragraphs for distinctive text\n# (robust to paragraph index drift) and then using a paragraph-scoped\n# Find/Replace so only the intended occurrence is touched:\n#   1. Professional summary sentence -> plain text swap to \"50M voters,\"\n#   2. Bulleted achievement -> \"50M\" becomes its own bold/colored run\n#   3. Impact statement -> plain text swap to \"50M voters nationwide,\"\n\n$d = $word.ActiveDocument\n\n# RGB() helper mirroring VBA's RGB(): Word's OLE_COLOR packs R + G*256 + B*65536.\nfunction RGB([int]$r, [int]$g, [int]$b) {\n    return $r + ($g * 256) + ($b * 65536)\n}\n\n$oldPhrase = \"all Black and Asian-American voters\"\n\nforeach ($p in $d.Paragraphs) {\n    $t = $p.Range.Text\n\n    if (-not $t.Contains($oldPhrase)) {\n        continue\n    }\n\n    if ($t.Contains(\"Discovered systematic race coding errors\")) {\n        # --- Edit 2: bulleted achievement --------------------------------\n        # Replace the whole phrase with plain neutral text first, then\n        # re-find just \"50M\" inside this paragraph and give it its own\n        # bold, colored run (matching the styling already used for the\n        # \"23%\"/\"64%\" runs later in the same bullet).\n        $rng = $p.Range\n        $null = $rng.Find.Execute($oldPhrase, $false, $false, $false, $false, $false, $true, 1, $false, \"50M voters\", 2)\n\n        $rng2 = $p.Range\n        $null = $rng2.Find.Execute(\"50M\", $false, $false, $false, $false, $false, $true, 1, $false, \"\", 0)\n        if ($rng2.Find.Found) {\n            $rng2.Font.Bold = 1\n            $rng2.Font.Color = RGB 0x2C 0x3E 0x50\n        }\n    }\n    elseif ($t.Contains(\"Product-focused\")) {\n        # --- Edit 1: professional summary sentence -----------------------\n        $rng = $p.Range\n        $null = $rng.Find.Execute(\"$oldPhrase,\", $false, $false, $false, $false, $false, $true, 1, $false, \"50M voters,\", 2)\n    }\n    elseif ($t.Contains(\"Impact:\")) {\n        # --- Edit 3: impact statement --------------------------------------\n        $rng = $p.Range\n        $null = $rng.Find.Execute(\"$oldPhrase,\", $false, $false, $false, $false, $false, $true, 1, $false, \"50M voters nationwide,\", 2)\n    }\n}\n"}
